# TC05_C3DC_phs002431 test case update:
#  1) Simplify the "Treatment Agent" SQL expression in the TreatmentTab query
#     (cell B5) by dropping the redundant outer CONCAT() wrapper around
#     REPLACE(...).
#  2) Re-apply the "wrap text" formatting to the other query cells
#     (B2, B3, B4, B6, B7) so they pick up the same cell style as B5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Fix the Treatment Agent SQL snippet inside the Treatment query text.
$treatmentQueryCell = $ws.Range("B5")
$sql = $treatmentQueryCell.Value2
$oldSnippet = "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))"
$newSnippet = "REPLACE(trt.treatment_agent, ';', ', ')"
$treatmentQueryCell.Value2 = $sql.Replace($oldSnippet, $newSnippet)

# 2) Refresh formatting on the other query cells so their style matches B5.
$queryCells = @("B2", "B3", "B4", "B6", "B7")
foreach ($addr in $queryCells) {
    $ws.Range($addr).WrapText = $true
}
